# Update the "Training Dashboard" sheet with the new progress as of 04-Nov-2025
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Make sure the "LAST UPDATE" column keeps storing plain text dates,
# not Excel date serials, to match the original inline-string format.
$ws.Range("I3:I5").NumberFormat = "@"

# Row 3
$ws.Range("H3").Value = 469
$ws.Range("I3").Value = "04-Nov-2025"

# Row 4
$ws.Range("H4").Value = -104
$ws.Range("I4").Value = "04-Nov-2025"

# Row 5
$ws.Range("H5").Value = 699
$ws.Range("I5").Value = "04-Nov-2025"
